# Features_Test.xlsx - "Keypoints deteketieren (SIFT) hinzugefügt"
#
# Inserts a new column E ("Anzahl Keypoints") before the existing "Labels"
# column (which shifts from E to F), and fills the new column with the
# SIFT keypoint counts for each of the 99 data rows (rows 2..100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing column E ("Labels" + its data) one column to the
# right, making room for the new "Anzahl Keypoints" column at E.
$ws.Columns.Item(5).Insert()

# New column header
$ws.Cells.Item(1, 5).Value = "Anzahl Keypoints"

# Keypoint counts for rows 2..100 (data rows 1..99)
$keypoints = @(197,204,138,218,121,279,184,144,294,225,220,271,235,243,122,118,163,106,95,260,224,183,174,184,131,289,220,230,234,236,244,125,93,360,260,158,130,103,118,126,151,367,266,188,141,111,131,98,138,154,224,273,301,239,216,154,107,219,151,151,299,250,214,132,128,152,240,321,229,248,213,262,258,239,125,96,159,136,94,84,108,182,77,145,139,112,148,111,74,58,26,117,38,151,93,76,75,87,157)

for ($i = 0; $i -lt $keypoints.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $keypoints[$i]
}
